# Update the "dSF" column (column F) values on Sheet1 to reflect the
# repulled / recalculated data, per the commit: "repull data, push all
# data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -3
    3  = 5
    4  = -5
    5  = 6
    6  = -1
    7  = -4
    8  = 1
    9  = -3
    10 = -5
    11 = 2
    13 = -2
    16 = -3
    17 = -1
    18 = 4
    19 = -5
    20 = 3
    22 = -4
    23 = 2
    24 = -1
    25 = 1
    27 = -4
    28 = -3
    29 = 3
    30 = -1
    32 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
